# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header in G1) holds per-game strikeout counts for
# this pitcher's 2023 save_data log. The values were regenerated from the
# authoritative box-score source (Strike# -> K), so we overwrite each game
# row's G value with the newly computed number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..68 (one entry per game row, in sheet order).
$kValues = @(
    2, 1, 1, 0, 0, 3, 1, 1, 2, 2,
    3, 0, 0, 1, 2, 1, 3, 1, 1, 1,
    0, 0, 2, 1, 1, 0, 0, 2, 2, 0,
    3, 1, 2, 1, 0, 2, 1, 0, 1, 1,
    1, 1, 0, 1, 0, 0, 2, 3, 1, 0,
    0, 1, 1, 1, 1, 4, 1, 2, 0, 1,
    1, 1, 1, 1, 2, 1, 0
)

$startRow = 2
$col = 7  # column G

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, $col).Value = $kValues[$i]
}
